$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row (row 4) that completes data collection
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = 16
$ws.Range("C4").Value = 4

# Update the active selection to reflect the next empty row, as Excel would
# after entering data and pressing Enter down the column (matches diff: B5)
$ws.Range("B5").Select()
